$wb = $excel.ActiveWorkbook

# Row 18 on ALC (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 250000510
$ws.Range("I18").Value = 690.6667
$ws.Range("K18").Value = 690.6667
$ws.Range("M18").Value = -406.6667

# Row 103 on ALC (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 5954135
$ws.Range("I103").Value = 2624.2856
$ws.Range("J103").Value = 14286251
$ws.Range("K103").Value = 7872.8568
$ws.Range("L103").Value = 42858753
$ws.Range("M103").Value = -7286.8568
$ws.Range("N103").Value = -42859925

# Row 113 on ALC (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3585.1667
$ws.Range("I113").Value = 2357
$ws.Range("K113").Value = 2357
$ws.Range("M113").Value = 897

# Row 116 on ALC (hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4481.364
$ws.Range("J116").Value = 4223.364
$ws.Range("L116").Value = 4223.364
$ws.Range("N116").Value = -11107.364

# Row 135 on ALC (hunk 4)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 682.2
$ws.Range("I135").Value = 466
$ws.Range("K135").Value = 4194
$ws.Range("M135").Value = -1659

# Row 141 on ALC (hunk 5)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2028.6786
$ws.Range("I141").Value = 1904.5416
$ws.Range("K141").Value = 5713.6248
$ws.Range("M141").Value = -533.6247999999996

# Row 45 on ARM (hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 39844.273
$ws.Range("J45").Value = 1813.8
$ws.Range("L45").Value = 1813.8
$ws.Range("N45").Value = -2567.8

# Row 74 on ARM (hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2178.946
$ws.Range("I74").Value = 1913.129
$ws.Range("J74").Value = 3552.3333
$ws.Range("K74").Value = 1913.129
$ws.Range("L74").Value = 3552.3333
$ws.Range("M74").Value = -1039.129
$ws.Range("N74").Value = -5300.3333

# Row 77 on ARM (hunk 8)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2178.946
$ws.Range("I77").Value = 1913.129
$ws.Range("J77").Value = 3552.3333
$ws.Range("K77").Value = 9565.645
$ws.Range("L77").Value = 17761.6665
$ws.Range("M77").Value = -5197.645
$ws.Range("N77").Value = -26497.6665

# Row 122 on ARM (hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3315.25
$ws.Range("I122").Value = 3315.25
$ws.Range("K122").Value = 9945.75
$ws.Range("M122").Value = -7495.75

# Row 132 on ARM (hunk 10)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1672.2858
$ws.Range("I132").Value = 980.3
$ws.Range("J132").Value = 2594.9333
$ws.Range("K132").Value = 2940.9
$ws.Range("L132").Value = 7784.7999
$ws.Range("M132").Value = -410.8999999999996
$ws.Range("N132").Value = -12844.7999

# Row 86 on BSM (hunk 11)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4371.857
$ws.Range("I86").Value = 3870.6
$ws.Range("K86").Value = 3870.6
$ws.Range("M86").Value = -2747.6

# Row 89 on BSM (hunk 12)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4371.857
$ws.Range("I89").Value = 3870.6
$ws.Range("K89").Value = 19353
$ws.Range("M89").Value = -13737

# Row 105 on BSM (hunk 13)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 21668598
$ws.Range("I105").Value = 1430316.9
$ws.Range("K105").Value = 1430316.9
$ws.Range("M105").Value = -1428569.9

# Row 107 on BSM (hunk 14)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2653859.2
$ws.Range("I107").Value = 4525971.5
$ws.Range("J107").Value = 1700.0834
$ws.Range("K107").Value = 4525971.5
$ws.Range("L107").Value = 1700.0834
$ws.Range("M107").Value = -4524051.5
$ws.Range("N107").Value = -5540.0834

# Row 31 on CRP (hunk 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5339.6855
$ws.Range("I31").Value = 3656.6
$ws.Range("K31").Value = 3656.6
$ws.Range("M31").Value = -3361.6

# Row 34 on CRP (hunk 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5339.6855
$ws.Range("I34").Value = 3656.6
$ws.Range("K34").Value = 3656.6
$ws.Range("M34").Value = -3454.6

# Row 58 on CRP (hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2088.889
$ws.Range("I58").Value = 1260.9166
$ws.Range("K58").Value = 1260.9166
$ws.Range("M58").Value = -1057.9166

# Row 62 on CRP (hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2746.1667
$ws.Range("I62").Value = 2993
$ws.Range("J62").Value = 2696.8
$ws.Range("K62").Value = 2993
$ws.Range("L62").Value = 2696.8
$ws.Range("M62").Value = -2369
$ws.Range("N62").Value = -3944.8

# Row 65 on CRP (hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2746.1667
$ws.Range("I65").Value = 2993
$ws.Range("J65").Value = 2696.8
$ws.Range("K65").Value = 14965
$ws.Range("L65").Value = 13484
$ws.Range("M65").Value = -11845
$ws.Range("N65").Value = -19724

# Row 122 on CRP (hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2816.75
$ws.Range("I122").Value = 3089
$ws.Range("K122").Value = 9267
$ws.Range("M122").Value = -6817

# Row 134 on CRP (hunk 21)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 11624.25
$ws.Range("I134").Value = 10499
$ws.Range("K134").Value = 31497
$ws.Range("M134").Value = -28962

# Row 136 on CRP (hunk 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2088.889
$ws.Range("I136").Value = 1260.9166
$ws.Range("K136").Value = 3782.7498
$ws.Range("M136").Value = -1232.7498

# Row 60 on CUL (hunk 23)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1541660.6
$ws.Range("I60").Value = 6667133.5
$ws.Range("K60").Value = 20001400.5
$ws.Range("M60").Value = -20001149.5

# Row 108 on CUL (hunk 24)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 686
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 121 on CUL (hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 4596431.5
$ws.Range("J121").Value = 65677.64999999999
$ws.Range("L121").Value = 197032.95
$ws.Range("N121").Value = -199652.95

# Row 131 on CUL (hunk 26)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 9266368
$ws.Range("J131").Value = 2152.7646
$ws.Range("L131").Value = 6458.293799999999
$ws.Range("N131").Value = -16538.2938

# Row 24 on GSM (hunk 27)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 6006
$ws.Range("I24").Value = 6006
$ws.Range("K24").Value = 6006
$ws.Range("M24").Value = -5833

# Row 70 on GSM (hunk 28)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 31255456
$ws.Range("I70").Value = 71432850
$ws.Range("J70").Value = 6376.8887
$ws.Range("K70").Value = 71432850
$ws.Range("L70").Value = 6376.8887
$ws.Range("M70").Value = -71432580
$ws.Range("N70").Value = -6916.8887

# Row 73 on GSM (hunk 29)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 31255456
$ws.Range("I73").Value = 71432850
$ws.Range("J73").Value = 6376.8887
$ws.Range("K73").Value = 71432850
$ws.Range("L73").Value = 6376.8887
$ws.Range("M73").Value = -71431914
$ws.Range("N73").Value = -8248.8887

# Row 113 on GSM (hunk 30)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1869.6
$ws.Range("J113").Value = 1706.5
$ws.Range("L113").Value = 1706.5
$ws.Range("N113").Value = -6046.5

# Row 16 on LTW (hunk 31)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1745.9
$ws.Range("I16").Value = 784.2857
$ws.Range("K16").Value = 784.2857
$ws.Range("M16").Value = -614.2857

# Row 93 on LTW (hunk 32)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1363
$ws.Range("I93").Value = 203
$ws.Range("K93").Value = 203
$ws.Range("M93").Value = 1045

# Row 132 on LTW (hunk 33)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9232
$ws.Range("I132").Value = 5174.2
$ws.Range("J132").Value = 13289.8
$ws.Range("K132").Value = 15522.6
$ws.Range("L132").Value = 39869.39999999999
$ws.Range("M132").Value = -12992.6
$ws.Range("N132").Value = -44929.39999999999

# Row 136 on LTW (hunk 34)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8200.666999999999
$ws.Range("I136").Value = 5961.6
$ws.Range("J136").Value = 10999.5
$ws.Range("K136").Value = 17884.8
$ws.Range("L136").Value = 32998.5
$ws.Range("M136").Value = -15334.8
$ws.Range("N136").Value = -38098.5

# Row 15 on WVR (hunk 35)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

# Row 63 on WVR (hunk 36)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 37500
$ws.Range("J63").Value = 37500
$ws.Range("L63").Value = 37500
$ws.Range("N63").Value = -38748

# Row 66 on WVR (hunk 37)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 37500
$ws.Range("J66").Value = 37500
$ws.Range("L66").Value = 112500
$ws.Range("N66").Value = -118740

# Row 69 on WVR (hunk 38)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 9750
$ws.Range("J69").Value = 9750
$ws.Range("L69").Value = 9750
$ws.Range("N69").Value = -11248

# Row 72 on WVR (hunk 39)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 9750
$ws.Range("J72").Value = 9750
$ws.Range("L72").Value = 29250
$ws.Range("N72").Value = -36738

# Row 132 on WVR (hunk 40)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2166.5173
$ws.Range("I132").Value = 2068.4
$ws.Range("J132").Value = 2779.75
$ws.Range("K132").Value = 6205.200000000001
$ws.Range("L132").Value = 8339.25
$ws.Range("M132").Value = -3675.200000000001
$ws.Range("N132").Value = -13399.25
